# "atualizacoes e pequenas correcoes"
#
#  1. Removes the _GoBack bookmark from the first paragraph ("Criar o
#     arquivo insperder ..."), leaving its text/proofErr runs intact.
#  2. Inserts a new block of paragraphs (a short .json save/load code
#     snippet) right after the 2nd of the three originally-empty
#     paragraphs that follow "Cria funcao que mostra todos os dados
#     salvos ate agora", i.e. before the 3rd (last) empty paragraph
#     that precedes "Importa arquivo insperdex...".
#  3. Re-creates the _GoBack bookmark, now living alone on its own
#     paragraph near the end of that new block.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Strip the bookmark out of paragraph 1, keeping its text/proofErr
#    runs exactly as they were.
# ---------------------------------------------------------------------
$p1xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Criar o arquivo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>insperder</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$d.Paragraphs.Item(1).Range.InsertXML($p1xml) | Out-Null

# ---------------------------------------------------------------------
# 2) Open up 7 fresh empty paragraphs right after paragraph 5 (the 2nd
#    empty paragraph following "Cria funcao..."). Character offsets
#    (Range.Start/End) are used instead of Paragraph.Index, which this
#    host does not maintain reliably across edits.
# ---------------------------------------------------------------------
$pos = $d.Paragraphs.Item(5).Range.End
for ($i = 0; $i -lt 7; $i++) {
    $r = $d.Range($pos, $pos)
    $r.InsertParagraphAfter() | Out-Null
    $pos = $pos + 1
}

# The newly-opened paragraphs are now #6-#12; fill the first four with
# the code snippet, leave #10-#11 empty, and put the relocated
# bookmark on #12. Paragraph #13 is the original 3rd empty paragraph,
# still sitting right before "Importa arquivo insperdex...".
$codeXml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Transforma </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>string</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>em</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

$codeXml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>arq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> = open("</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>teste.json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>", "w")</w:t></w:r></w:p>
'@

$codeXml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>arq.write</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>meunome</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@

$codeXml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>arq.close</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>()</w:t></w:r></w:p>
'@

$bookmarkXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

# A freshly InsertParagraphAfter-ed paragraph round-trips as <w:p><w:r/></w:p>;
# normalize the ones that should stay genuinely empty back to a bare <w:p/>
# (what Word itself emits for an empty paragraph), matching #10/#11/#13.
$emptyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$d.Paragraphs.Item(6).Range.InsertXML($codeXml1) | Out-Null
$d.Paragraphs.Item(7).Range.InsertXML($codeXml2) | Out-Null
$d.Paragraphs.Item(8).Range.InsertXML($codeXml3) | Out-Null
$d.Paragraphs.Item(9).Range.InsertXML($codeXml4) | Out-Null
$d.Paragraphs.Item(10).Range.InsertXML($emptyXml) | Out-Null
$d.Paragraphs.Item(11).Range.InsertXML($emptyXml) | Out-Null
$d.Paragraphs.Item(12).Range.InsertXML($bookmarkXml) | Out-Null

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
